# Estadisticos Matutinos 15 Oct
# Add a rescued-student record to the "Rescatables" sheet (4th worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New row 2: NC | Paterno | Materno | Nombres | Nombre_Largo | Grupo | Reprobadas
$ws.Cells.Item(2, 1).Value = 20330051920112
$ws.Cells.Item(2, 2).Value = "VASQUEZ"
$ws.Cells.Item(2, 3).Value = "PEREZ"
$ws.Cells.Item(2, 4).Value = "GUILLERMO"
$ws.Cells.Item(2, 5).Value = "MANTIENE LOS MOTORES DE CA Y CC"
$ws.Cells.Item(2, 6).Value = "3BEM"
$ws.Cells.Item(2, 7).Value = 6
